$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove columns D and E (promocionado_25, Black_Friday)
$ws.Range("D1:E9").EntireColumn.Delete()

# Update Preco (C) column values
$ws.Range("C2").Value = 379
$ws.Range("C4").Value = 370
$ws.Range("C5").Value = 369
$ws.Range("C6").Value = 329
$ws.Range("C7").Value = 372
$ws.Range("C8").Value = 369
$ws.Range("C9").Value = 366

# Add new rows 10 and 11
$ws.Range("A10").Value = 45926
$ws.Range("B10").Value = 7172
$ws.Range("C10").Value = 358

$ws.Range("A11").Value = 45927
$ws.Range("B11").Value = 7172
$ws.Range("C11").Value = 358

# Apply same date style (style index 1) to A10:A11, copy format from A9
$ws.Range("A9").Copy()
$ws.Range("A10:A11").PasteSpecial(-4122) | Out-Null # xlPasteFormats

# Set the active selection
$ws.Range("D4").Select() | Out-Null
